$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 16-23 (no longer present in the updated sheet)
$ws.Range("A16:F23").EntireRow.Delete()

# Row 2
$ws.Range("B2").Value = "NSE:ATFL"
$ws.Range("C2").Value = "NSE:BEDMUTHA"
$ws.Range("E2").Value = "NSE:MOTHERSON"
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("B3").Value = "NSE:INDSWFTLAB"
$ws.Range("C3").Value = "NSE:CANTABIL"
$ws.Range("E3").Value = ""

# Row 4
$ws.Range("B4").Value = "NSE:NITIRAJ"
$ws.Range("C4").Value = "NSE:CREATIVEYE"

# Row 5
$ws.Range("C5").Value = "NSE:EKC"

# Row 6
$ws.Range("C6").Value = "NSE:HBSL"

# Row 7
$ws.Range("C7").Value = "NSE:IMAGICAA"

# Row 8
$ws.Range("C8").Value = "NSE:INDORAMA"

# Row 9
$ws.Range("C9").Value = "NSE:KILITCH"

# Row 10
$ws.Range("C10").Value = "NSE:OILCOUNTUB"

# Row 11
$ws.Range("C11").Value = "NSE:PAGEIND"

# Row 12
$ws.Range("C12").Value = "NSE:PRINCEPIPE"

# Row 13
$ws.Range("C13").Value = "NSE:RAJSREESUG"

# Row 14
$ws.Range("C14").Value = "NSE:RAMAPHO"

# Row 15
$ws.Range("C15").Value = "NSE:SALSTEEL"
